# Horarios actualizados Linea 141 - 311
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Update header info
$ws1.Range("A2").Value = "Última actualización: 05:18:42"
$ws1.Range("A3").Value = "Total filas: 32"

# Insert a new row at position 18 (pushes old rows 18-31 down to 19-32)
$ws1.Rows.Item(18).Insert()

# Fill the newly inserted row 18 with the new scraped entry
$ws1.Cells.Item(18, 1).Value = "05:18:42"
$ws1.Cells.Item(18, 2).Value = "05:25"
$ws1.Cells.Item(18, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(18, 4).Value = 7
$ws1.Cells.Item(18, 5).Value = "LP1912"

# Append 5 new rows at the bottom (rows 33-37)
$ws1.Cells.Item(33, 1).Value = "05:18:42"
$ws1.Cells.Item(33, 2).Value = "06:59"
$ws1.Cells.Item(33, 3).Value = "14_ABASTO"
$ws1.Cells.Item(33, 4).Value = 101
$ws1.Cells.Item(33, 5).Value = "LP1912"

$ws1.Cells.Item(34, 1).Value = "05:18:42"
$ws1.Cells.Item(34, 2).Value = "07:05"
$ws1.Cells.Item(34, 3).Value = "15_ABASTO"
$ws1.Cells.Item(34, 4).Value = 107
$ws1.Cells.Item(34, 5).Value = "LP1912"

$ws1.Cells.Item(35, 1).Value = "05:18:42"
$ws1.Cells.Item(35, 2).Value = "07:07"
$ws1.Cells.Item(35, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(35, 4).Value = 109
$ws1.Cells.Item(35, 5).Value = "LP1912"

$ws1.Cells.Item(36, 1).Value = "05:18:42"
$ws1.Cells.Item(36, 2).Value = "07:11"
$ws1.Cells.Item(36, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(36, 4).Value = 113
$ws1.Cells.Item(36, 5).Value = "LP1912"

$ws1.Cells.Item(37, 1).Value = "05:18:42"
$ws1.Cells.Item(37, 2).Value = "07:15"
$ws1.Cells.Item(37, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(37, 4).Value = 117
$ws1.Cells.Item(37, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

# Update header info
$ws2.Range("A2").Value = "Última actualización: 05:18:42"
$ws2.Range("A3").Value = "Total filas: 11"

# Append 1 new row at the bottom (row 16)
$ws2.Cells.Item(16, 1).Value = "05:18:42"
$ws2.Cells.Item(16, 2).Value = "07:11"
$ws2.Cells.Item(16, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(16, 4).Value = 113
$ws2.Cells.Item(16, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

# Update header info
$ws3.Range("A2").Value = "Última actualización: 05:18:42"
$ws3.Range("A3").Value = "Total filas: 8"

# Append 1 new row at the bottom (row 13)
$ws3.Cells.Item(13, 1).Value = "05:18:42"
$ws3.Cells.Item(13, 2).Value = "06:59"
$ws3.Cells.Item(13, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(13, 4).Value = 101
$ws3.Cells.Item(13, 5).Value = "L6173"
